# Insert a new "label" row into the key/value config sheet, between the
# "sig" row and the "headline" row (i.e. as the new row 3), pushing every
# row from the old row 3 onward down by one. This mirrors adding a new
# "label" field to the template (commit: "added label div ...").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 3 - shifts existing rows 3..12 down to 4..13
# and inherits formatting from the surrounding rows.
$ws.Rows("3:3").Insert()

# Populate the newly inserted row with the new key/value pair.
$ws.Range("A3").Value = "label"
$ws.Range("B3").Value = "Label goes here"
